$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 300 (pushes old rows 300-357 down to 302-359)
$ws.Range("A300:R301").EntireRow.Insert()

# Fill new row 300 with the new weekly record
$ws.Cells.Item(300, 1).Value = 5
$ws.Cells.Item(300, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(300, 3).Value = "Maule"
$ws.Cells.Item(300, 4).Value = 44476
$ws.Cells.Item(300, 5).Value = 7
$ws.Cells.Item(300, 6).Value = 100112020
$ws.Cells.Item(300, 7).Value = "Tomate"
$ws.Cells.Item(300, 8).Value = "Larga vida"
$ws.Cells.Item(300, 9).Value = "Primera"
$ws.Cells.Item(300, 10).Value = 2000
$ws.Cells.Item(300, 11).Value = 16000
$ws.Cells.Item(300, 12).Value = 16000
$ws.Cells.Item(300, 13).Value = 16000
$ws.Cells.Item(300, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(300, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(300, 16).Value = 889
$ws.Cells.Item(300, 17).Value = 18
$ws.Cells.Item(300, 18).Value = "Hortaliza"

# Fill new row 301 with the new weekly record
$ws.Cells.Item(301, 1).Value = 5
$ws.Cells.Item(301, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(301, 3).Value = "Maule"
$ws.Cells.Item(301, 4).Value = 44476
$ws.Cells.Item(301, 5).Value = 7
$ws.Cells.Item(301, 6).Value = 100112020
$ws.Cells.Item(301, 7).Value = "Tomate"
$ws.Cells.Item(301, 8).Value = "Larga vida"
$ws.Cells.Item(301, 9).Value = "Primera"
$ws.Cells.Item(301, 10).Value = 1500
$ws.Cells.Item(301, 11).Value = 7000
$ws.Cells.Item(301, 12).Value = 7000
$ws.Cells.Item(301, 13).Value = 7000
$ws.Cells.Item(301, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(301, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(301, 16).Value = 700
$ws.Cells.Item(301, 17).Value = 10
$ws.Cells.Item(301, 18).Value = "Hortaliza"

Write-Output "done"
